$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B ("NISN") to make room for "Jurusan"
$ws.Range("B:B").EntireColumn.Insert()

# --- New header cell (row 1, title row style: bold, size 16) ---
$ws.Range("B1").Font.Bold = $true
$ws.Range("B1").Font.Size = 16

# --- New "Jurusan" header cell (row 2, bold header style with yellow fill + border) ---
$ws.Range("B2").Value = "Jurusan"
$ws.Range("B2").Font.Bold = $true
$ws.Range("B2").Interior.Color = 65535
$ws.Range("B2").Borders.LineStyle = 1

# --- New data cells (rows 3-7, bordered style) ---
$ws.Range("B3").Value = "IPA"
$ws.Range("B4").Value = "IPA"
$ws.Range("B5").Value = "IPA"
$ws.Range("B6").Value = "IPS"
$ws.Range("B7").Value = "IPS"
$ws.Range("B3:B7").Borders.LineStyle = 1

# Restore a plausible final selection/cursor position (cosmetic, matches author's last click)
$null = $ws.Range("D16").Select()
